# Updates the "Price" (column D) and a couple "Volume(1h)" (column E) text
# values in the crypto symbol list sheet, refreshed by the scheduled
# GitHub Actions scraper run. Values are plain text (not numbers) in the
# source sheet, so NumberFormat is forced to "@" (Text) before assigning
# the string -- otherwise Excel would coerce a numeric-looking string like
# "244.30" into the number 244.3 and drop the significant trailing zero.
# The format is reset back to the default "Normal" style afterward so we
# don't leave a stray text-format style applied to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "244.30"
$ws.Range("D2").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.385"
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06003"
$ws.Range("D5").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8186"
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9541"
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1424"
$ws.Range("D9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07438"
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03283"
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03044"
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09407"
$ws.Range("D13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.005"
$ws.Range("D14").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04823"
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0005912"
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.005450"
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.004151"
$ws.Range("D19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0009858"
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.673"
$ws.Range("D21").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.416"
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.190"
$ws.Range("D23").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03995"
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006422"
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1072"
$ws.Range("D42").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.005815"
$ws.Range("D44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005108"
$ws.Range("D45").Style = "Normal"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.004614"
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("D49").Style = "Normal"

